# Append 5 new abituryent (applicant) rows (233-237) to the sheet,
# matching the source data exactly. Columns E, F, I, J, K contain
# numeric/date-looking text (passport series+number, JSHIR id, phone
# numbers, ISO date) that must stay as literal text, not be coerced
# into numbers/dates by Excel, so we force NumberFormat "@" on those
# columns before assigning values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ r = 233; A = "Ahmadjonov  Doniyorbek Rustamovich";  B = "Yurisprudensiya";    C = "Rus tili";    D = "Kunduzgi"; E = "AD0562512"; F = "52107050005111"; G = "Samarqand viloyati"; H = "Samarqand tumani"; I = "998992240337"; J = "+998992240337"; K = "2025-07-21" },
    @{ r = 234; A = "Abduraimov Shohjahon Begzodivoch";    B = "Yurisprudensiya";    C = "O'zbek tili"; D = "Kunduzgi"; E = "AE1322337"; F = "50109085540038"; G = "Toshkent shahri";     H = "Mirzo Ulugʻbek tumani"; I = "998911997555"; J = "+998958127100"; K = "2025-07-21" },
    @{ r = 235; A = "Turg'unboyev Shohrux Sardorovich";    B = "Bugalteriya hisobi"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AD6720333"; F = "52011076830021"; G = "Toshkent shahri";     H = "Chilonzor tumani"; I = "998908270020"; J = "+998908270020"; K = "2025-07-21" },
    @{ r = 236; A = "Abduraxmanov ibroximjon rustamovich"; B = "Yurisprudensiya";    C = "Rus tili";    D = "Kunduzgi"; E = "AD0564279"; F = "31305795040012"; G = "Andijon viloyati";    H = "Andijon tuman";    I = "998979996656"; J = "+998502776657"; K = "2025-07-21" },
    @{ r = 237; A = "Yulduzova Farida";                    B = "Yurisprudensiya";    C = "O'zbek tili"; D = "Kunduzgi"; E = "AD7674668"; F = "62107075330022"; G = "Buxoro viloyati";      H = "Peshku tumani";    I = "998993842777"; J = "+998992961425"; K = "2025-07-22" }
)

foreach ($row in $rows) {
    $r = $row.r

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D

    # Passport series/number (e.g. AD0562512) -> keep as text
    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = $row.E

    # JSHIR / PINFL (14-digit id) -> keep as text, not a number
    $ws.Cells.Item($r, 6).NumberFormat = "@"
    $ws.Cells.Item($r, 6).Value = $row.F

    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H

    # Telegram number -> keep as text
    $ws.Cells.Item($r, 9).NumberFormat = "@"
    $ws.Cells.Item($r, 9).Value = $row.I

    # Phone number (with leading +) -> keep as text
    $ws.Cells.Item($r, 10).NumberFormat = "@"
    $ws.Cells.Item($r, 10).Value = $row.J

    # Date string "YYYY-MM-DD" -> keep as literal text, not converted to a date serial
    $ws.Cells.Item($r, 11).NumberFormat = "@"
    $ws.Cells.Item($r, 11).Value = $row.K
}
